$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so formatting (trailing zeros, thousand dots, etc.) is preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.669.86'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '1.863.79'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '0.7236'
$ws.Range("E5").Value = '  -1.86%  '
$ws.Range("D6").Value = '239.69'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("D7").Value = '0.9990'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.3104'
$ws.Range("E8").Value = '  -1.56%  '
$ws.Range("D9").Value = '0.07042'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = '24.17'
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").Value = '0.08185'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").Value = '0.7374'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '5.276'
$ws.Range("E13").Value = '  -3.19%  '
$ws.Range("D14").Value = '1.843.34'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").Value = '91.60'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '29.662.89'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '5.968'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '245.54'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").Value = '13.27'
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").Value = '0.000007745'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").Value = '0.9987'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '2.101.38'
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").Value = '0.9990'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '7.671'
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = '0.1521'
$ws.Range("E25").Value = '  -2.70%  '
$ws.Range("D26").Value = '9.115'
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("D27").Value = '162.44'
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").Value = '18.38'
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").Value = '1.994'
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '1.434'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.517'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.475'
$ws.Range("E32").Value = '  -4.01%  '
$ws.Range("D33").Value = '4.148'
$ws.Range("E33").Value = '  -3.60%  '
$ws.Range("D34").Value = '0.05240'
$ws.Range("E34").Value = '  -1.66%  '
$ws.Range("D35").Value = '1.222'
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("D36").Value = '0.7451'
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").Value = '0.9972'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '2.691'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '0.01916'
$ws.Range("E39").Value = '  -1.72%  '
$ws.Range("D40").Value = '2.729'
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").Value = '0.4432'
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("D42").Value = '5.962'
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("D43").Value = '0.8622'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '70.54'
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("D45").Value = '1.044.68'
$ws.Range("E45").Value = '  -5.11%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9994'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '103.67'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("D48").Value = '7.411'
$ws.Range("E48").Value = '  -3.64%  '
$ws.Range("D49").Value = '1.804'
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.451'
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.001.04'
$ws.Range("E51").Value = '  -0.63%  '
